$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the combined resistor reference list (R2,R3,R9,R10,R11,R21,R22,R23,R27,R28)
# into two groups: the 0402 group (R21,R22,R23) moves into the previously-empty
# template row 8, and the remaining 0603 group (R2,R3,R9,R10,R11,R27,R28) stays on row 9.

# Row 8: fill in the reference designator and SMD parts count for the new 0402 group.
$ws.Range("A8").Value2 = "R21,R22,R23"
$ws.Range("J8").Value2 = 3

# Row 9: narrow the reference list to the remaining resistors and update the count.
$ws.Range("A9").Value2 = "R2,R3,R9,R10,R11,R27,R28"
$ws.Range("J9").Value2 = 7

# Row 9 no longer needs the taller, manually-set row height now that the text is shorter.
$ws.Rows("9").AutoFit()

# Update the active selection / scroll position recorded with the sheet view.
$ws.Range("J10").Select()
try {
    $win = $excel.ActiveWindow
    $win.ScrollColumn = 5
    $win.ScrollRow = 1
} catch {
}
